$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cobertura de la Prueba")

$ws.Range("D3").Value = 69

$ws.Range("D4").Select()
